$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells ---

# Row 3: "Temps [h]" was entered as a plain number (1); change it to the
# French-style text "1,00" like the later entries in the journal.
$ws.Range("C3").Value = "1,00"

# Row 12: hours corrected from "0,15" to "0,25" (description/type unchanged).
$ws.Range("C12").Value = "0,25"

# --- Append new journal entries (rows 14-21) ---

# Seed formatting for the new rows by copying the format of the last
# populated row (row 13), then overwrite the values.
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F21").PasteSpecial(-4122)

# Row 14
$ws.Range("A14").Value = 44977
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = "0,30"
$ws.Range("D14").Value = "Analyse"
$ws.Range("E14").Value = "Discution avec le maître de projet"

# Row 15
$ws.Range("A15").Value = 44977
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = "1,15"
$ws.Range("D15").Value = "Implémentation"
$ws.Range("E15").Value = "Commencer le front-end de la page d'accueil"

# Row 16
$ws.Range("A16").Value = 44977
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = "2,25"
$ws.Range("D16").Value = "Implémentation"
$ws.Range("E16").Value = "Remodéliser la base de donnée"

# Row 17
$ws.Range("A17").Value = 45006
$ws.Range("B17").Value = 3
$ws.Range("C17").Value = "1,50"
$ws.Range("D17").Value = "Implémentation"
$ws.Range("E17").Value = "Mettre à jour le MCD et MLD"

# Row 18
$ws.Range("A18").Value = 45008
$ws.Range("B18").Value = 3
$ws.Range("C18").Value = 1.5
$ws.Range("D18").Value = "Implémentation"
$ws.Range("E18").Value = "Création et modification des backlogs"

# Row 19
$ws.Range("A19").Value = 45009
$ws.Range("B19").Value = 3
$ws.Range("C19").Value = "0,75"
$ws.Range("D19").Value = "Implémentation"
$ws.Range("E19").Value = "Compléter les infomations du IceScrum"

# Row 20
$ws.Range("A20").Value = 45009
$ws.Range("B20").Value = 3
$ws.Range("C20").Value = "1,00"
$ws.Range("D20").Value = "Documentation"
$ws.Range("E20").Value = "Parfaire le rapport de projet"

# Row 21
$ws.Range("A21").Value = 45009
$ws.Range("B21").Value = 3
$ws.Range("C21").Value = "0,50"
$ws.Range("D21").Value = "Implémentation"
$ws.Range("E21").Value = "Modifier les maquettes"

# Match the saved selection from the source edit (last cell touched).
$ws.Range("E21").Select()
